# Weekly "Förändrad" (changed) date bump + one new logging-notice row.
#
# The report's "Förändrad" column (C) is stamped with the date the report
# was regenerated; every existing record (rows 2-305) moves from
# 2023-09-11 (45180) to 2023-09-12 (45181). A brand-new notice
# (A 42336-2023) is appended as row 306.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" date for every existing data row.
$ws.Range("C2:C305").Value = 45181

# 2) Row 305 picks up an explicit row height once the sheet grows past it
#    (matches the height already used by every other data row).
$ws.Rows.Item(305).RowHeight = 15

# 3) Append the new logging notice as row 306.
$ws.Cells.Item(306, 1).Value = "A 42336-2023"

$ws.Cells.Item(306, 2).Value = 45180
$ws.Cells.Item(306, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(306, 3).Value = 45181
$ws.Cells.Item(306, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(306, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(306, 5).Value = "SÖDERHAMN"
$ws.Cells.Item(306, 6).Value = "Bergvik skog väst AB"

$ws.Cells.Item(306, 7).Value = 3
$ws.Cells.Item(306, 8).Value = 0
$ws.Cells.Item(306, 9).Value = 0
$ws.Cells.Item(306, 10).Value = 0
$ws.Cells.Item(306, 11).Value = 0
$ws.Cells.Item(306, 12).Value = 0
$ws.Cells.Item(306, 13).Value = 0
$ws.Cells.Item(306, 14).Value = 0
$ws.Cells.Item(306, 15).Value = 0
$ws.Cells.Item(306, 16).Value = 0
$ws.Cells.Item(306, 17).Value = 0

# R306 ("Artnamn") stays blank but keeps the wrap-text style used by the
# rest of the column.
$ws.Cells.Item(306, 18).WrapText = $true
$ws.Cells.Item(306, 18).Value = ""
